# Edit script: add new column M "OV_UK_VISITORS_NSA" (UK overseas visitors, NSA)
# to the "Data" worksheet, refresh column widths / selection / window view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. New header cell M1, sharing the same header style as the rest of row 1 ---
$ws.Range("M1").Value = "OV_UK_VISITORS_NSA"
$ws.Range("L1").Copy()
$ws.Range("M1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- 2. New column M data, rows 2-199 ---
# Rows 2-41 are #N/A (no overseas-visitor series available that far back);
# rows 42-199 carry the actual OV_UK_VISITORS_NSA figures.
$mValues = @("#N/A","#N/A","#N/A","#N/A","#N/A","#N/A","#N/A","#N/A","#N/A","#N/A","#N/A","#N/A","#N/A","#N/A","#N/A","#N/A","#N/A","#N/A","#N/A","#N/A","#N/A","#N/A","#N/A","#N/A","#N/A","#N/A","#N/A","#N/A","#N/A","#N/A","#N/A","#N/A","#N/A","#N/A","#N/A","#N/A","#N/A","#N/A","#N/A","#N/A",2081,3240,4738,2360,1920,3008,4261,2262,2013,3174,4255,2196,2013,3200,4714,2537,2153,3582,5180,2727,2336,3958,5405,2751,2579,3320,5065,2933,2641,4047,5618,3259,2777,4013,5549,3460,3337,4265,5961,3776,3323,4537,6306,3851,2836,4297,6012,3980,3345,4897,6188,4103,3611,5222,6700,4333,3840,5269,6891,4794,4217,5952,7747,5622,4719,6680,8130,5631,4940,6446,8168,5961,4804,6834,8028,6080,5046,6800,7912,5636,4993,6733,7943,5538,4862,6279,7098,4593,4558,6509,7619,5706,4926,6085,7716,6147,5451,7015,8412,6748,6171,7868,8858,7073,6351,8476,10296,7592,6737,8511,9564,7966,7195,8406,9358,6929,6213,7943,8779,6951,5904,7755,8985,7160,6113,8308,9172,7205,6249,8457,8786,7593,6332,8869,9645,7843,6810,9490,10030,8050,6949,9983,10499,8683,7551,9968,10659,9432,8307,10877,11110,8920,7883,10056,10802,9163,7956,9593)

$mArray = New-Object 'object[,]' 198,1
for ($i = 0; $i -lt 198; $i++) {
    $mArray[$i,0] = $mValues[$i]
}
$ws.Range("M2:M199").Value = $mArray

# --- 3. Column widths (best-fit on the columns that matter) ---
$ws.Columns.Item(1).ColumnWidth = 9.25
$ws.Columns.Item(12).ColumnWidth = 26.59
$ws.Columns.Item(13).ColumnWidth = 17.75

# --- 4. Selection moves to H14 ---
$ws.Range("H14").Select()
